# Update the "dSF" column (F) values for specific rows, per repull/recalculation
# of the underlying data (commit: "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    6  = 0
    11 = -3
    14 = 1
    20 = 2
    22 = -1
    23 = 0
    24 = 2
    25 = 2
    28 = 3
    32 = 3
    33 = 1
    48 = 1
    51 = -4
    52 = 1
    57 = -6
    59 = -1
    60 = 0
    62 = 0
    64 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
